# Update the Monte-Carlo simulation results on the "results" sheet with the
# final, post-publication-quality-check numbers (columns C:O, rows 2-17).
# Column A/B (phack, gain) and the header row are unchanged; sheet
# "parameters" is unaffected by this refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("results")

$ws.Range("C2").Value = 0.19995522495688081
$ws.Range("D2").Value = 0.95119096111961443
$ws.Range("E2").Value = 9.0109999999999992
$ws.Range("F2").Value = 126.858
$ws.Range("G2").Value = 11.759
$ws.Range("H2").Value = 1.6619999999999999
$ws.Range("I2").Value = 1.1599999999999999
$ws.Range("J2").Value = 0.19845633406283003
$ws.Range("K2").Value = 0.94591973833270215
$ws.Range("L2").Value = 26.858000000000001
$ws.Range("M2").Value = 126.858
$ws.Range("N2").Value = 1.766
$ws.Range("O2").Value = 1.1850000000000001
$ws.Range("C3").Value = 0.20110348331762906
$ws.Range("D3").Value = 0.94311202997080368
$ws.Range("E3").Value = 8.9600000000000009
$ws.Range("F3").Value = 126.777
$ws.Range("G3").Value = 11.759
$ws.Range("H3").Value = 1.1739999999999999
$ws.Range("I3").Value = 0.80700000000000005
$ws.Range("J3").Value = 0.18141227793545692
$ws.Range("K3").Value = 0.94306811332127805
$ws.Range("L3").Value = 26.777000000000001
$ws.Range("M3").Value = 102.64700000000001
$ws.Range("N3").Value = 1.129
$ws.Range("O3").Value = 0.77100000000000002
$ws.Range("C4").Value = 0.20083990102901575
$ws.Range("D4").Value = 0.94680173188515138
$ws.Range("E4").Value = 9.0310000000000006
$ws.Range("F4").Value = 126.745
$ws.Range("G4").Value = 11.577
$ws.Range("H4").Value = 1.7549999999999999
$ws.Range("I4").Value = 1.1020000000000001
$ws.Range("J4").Value = 0.19788471530815782
$ws.Range("K4").Value = 0.95031085604973664
$ws.Range("L4").Value = 26.745000000000001
$ws.Range("M4").Value = 103.88200000000001
$ws.Range("N4").Value = 1.7230000000000001
$ws.Range("O4").Value = 1.18
$ws.Range("C5").Value = 0.2007045939949916
$ws.Range("D5").Value = 0.93816578711266607
$ws.Range("E5").Value = 9.0760000000000005
$ws.Range("F5").Value = 126.752
$ws.Range("G5").Value = 11.554
$ws.Range("H5").Value = 2.7149999999999999
$ws.Range("I5").Value = 1.89
$ws.Range("J5").Value = 0.1948407098003464
$ws.Range("K5").Value = 0.95027082609199376
$ws.Range("L5").Value = 26.751999999999999
$ws.Range("M5").Value = 106.258
$ws.Range("N5").Value = 2.7309999999999999
$ws.Range("O5").Value = 1.893
$ws.Range("C6").Value = 0.2490767973411048
$ws.Range("D6").Value = 0.96872069195599519
$ws.Range("E6").Value = 9.1010000000000009
$ws.Range("F6").Value = 131.446
$ws.Range("G6").Value = 14.157999999999999
$ws.Range("H6").Value = 2.0950000000000002
$ws.Range("I6").Value = 1.175
$ws.Range("J6").Value = 0.2471670456579611
$ws.Range("K6").Value = 0.97050664566628997
$ws.Range("L6").Value = 31.446000000000002
$ws.Range("M6").Value = 131.446
$ws.Range("N6").Value = 1.9530000000000001
$ws.Range("O6").Value = 1.226
$ws.Range("C7").Value = 0.24772486763301227
$ws.Range("D7").Value = 0.97475464792479416
$ws.Range("E7").Value = 8.9480000000000004
$ws.Range("F7").Value = 131.36099999999999
$ws.Range("G7").Value = 14.314
$ws.Range("H7").Value = 1.4350000000000001
$ws.Range("I7").Value = 0.88500000000000001
$ws.Range("J7").Value = 0.23075056202738609
$ws.Range("K7").Value = 0.9717314487632509
$ws.Range("L7").Value = 31.361000000000001
$ws.Range("M7").Value = 103.2
$ws.Range("N7").Value = 1.4470000000000001
$ws.Range("O7").Value = 0.878
$ws.Range("C8").Value = 0.25163833075460934
$ws.Range("D8").Value = 0.97126633543623475
$ws.Range("E8").Value = 8.91
$ws.Range("F8").Value = 131.53399999999999
$ws.Range("G8").Value = 14.090999999999999
$ws.Range("H8").Value = 2.1070000000000002
$ws.Range("I8").Value = 1.2150000000000001
$ws.Range("J8").Value = 0.24733980011757795
$ws.Range("K8").Value = 0.97337137422729403
$ws.Range("L8").Value = 31.533999999999999
$ws.Range("M8").Value = 104.59099999999999
$ws.Range("N8").Value = 2.0350000000000001
$ws.Range("O8").Value = 1.2190000000000001
$ws.Range("C9").Value = 0.24997677939366206
$ws.Range("D9").Value = 0.97079295249804265
$ws.Range("E9").Value = 8.9939999999999998
$ws.Range("F9").Value = 131.38399999999999
$ws.Range("G9").Value = 13.997
$ws.Range("H9").Value = 3.2
$ws.Range("I9").Value = 1.94
$ws.Range("J9").Value = 0.25006151209762018
$ws.Range("K9").Value = 0.96627124347972393
$ws.Range("L9").Value = 31.384
$ws.Range("M9").Value = 107.381
$ws.Range("N9").Value = 3.3180000000000001
$ws.Range("O9").Value = 1.9570000000000001
$ws.Range("C10").Value = 0.29620797865728682
$ws.Range("D10").Value = 0.99592380231514788
$ws.Range("E10").Value = 8.8819999999999997
$ws.Range("F10").Value = 135.74100000000001
$ws.Range("G10").Value = 16.558
$ws.Range("H10").Value = 2.4060000000000001
$ws.Range("I10").Value = 1.2769999999999999
$ws.Range("J10").Value = 0.29455206320537547
$ws.Range("K10").Value = 0.99429264393956263
$ws.Range("L10").Value = 35.741
$ws.Range("M10").Value = 135.74100000000001
$ws.Range("N10").Value = 2.54
$ws.Range("O10").Value = 1.248
$ws.Range("C11").Value = 0.29434673244948173
$ws.Range("D11").Value = 0.99566049317553107
$ws.Range("E11").Value = 8.8249999999999993
$ws.Range("F11").Value = 135.49199999999999
$ws.Range("G11").Value = 16.483000000000001
$ws.Range("H11").Value = 1.643
$ws.Range("I11").Value = 0.91600000000000004
$ws.Range("J11").Value = 0.27968036529680357
$ws.Range("K11").Value = 0.9929017603634297
$ws.Range("L11").Value = 35.491999999999997
$ws.Range("M11").Value = 103.702
$ws.Range("N11").Value = 1.706
$ws.Range("O11").Value = 0.92100000000000004
$ws.Range("C12").Value = 0.29552045439881047
$ws.Range("D12").Value = 0.99558387830188011
$ws.Range("E12").Value = 9.0489999999999995
$ws.Range("F12").Value = 135.827
$ws.Range("G12").Value = 16.78
$ws.Range("H12").Value = 2.5419999999999998
$ws.Range("I12").Value = 1.1619999999999999
$ws.Range("J12").Value = 0.29815491362747454
$ws.Range("K12").Value = 0.99883612662942267
$ws.Range("L12").Value = 35.826999999999998
$ws.Range("M12").Value = 105.161
$ws.Range("N12").Value = 2.4289999999999998
$ws.Range("O12").Value = 1.2909999999999999
$ws.Range("C13").Value = 0.29293136578730194
$ws.Range("D13").Value = 0.99690757525724827
$ws.Range("E13").Value = 9.1240000000000006
$ws.Range("F13").Value = 135.714
$ws.Range("G13").Value = 16.934999999999999
$ws.Range("H13").Value = 4.1669999999999998
$ws.Range("I13").Value = 2.137
$ws.Range("J13").Value = 0.28160348781263822
$ws.Range("K13").Value = 0.99825174825174823
$ws.Range("L13").Value = 35.713999999999999
$ws.Range("M13").Value = 108.42400000000001
$ws.Range("N13").Value = 3.8359999999999999
$ws.Range("O13").Value = 2.06
$ws.Range("C14").Value = 0.31820274499773671
$ws.Range("D14").Value = 0.99907144736804043
$ws.Range("E14").Value = 8.8290000000000006
$ws.Range("F14").Value = 137.71100000000001
$ws.Range("G14").Value = 17.997
$ws.Range("H14").Value = 2.472
$ws.Range("I14").Value = 1.306
$ws.Range("J14").Value = 0.321263597700434
$ws.Range("K14").Value = 0.99953958032615331
$ws.Range("L14").Value = 37.710999999999999
$ws.Range("M14").Value = 137.71100000000001
$ws.Range("N14").Value = 2.5779999999999998
$ws.Range("O14").Value = 1.206
$ws.Range("C15").Value = 0.31830708447355266
$ws.Range("D15").Value = 0.99936770893096205
$ws.Range("E15").Value = 9.0990000000000002
$ws.Range("F15").Value = 137.85
$ws.Range("G15").Value = 18.102
$ws.Range("H15").Value = 1.7729999999999999
$ws.Range("I15").Value = 0.89800000000000002
$ws.Range("J15").Value = 0.32492038870011947
$ws.Range("L15").Value = 37.85
$ws.Range("M15").Value = 103.77200000000001
$ws.Range("N15").Value = 1.8129999999999999
$ws.Range("O15").Value = 0.88800000000000001
$ws.Range("C16").Value = 0.31295049299874012
$ws.Range("D16").Value = 0.9997906197654941
$ws.Range("E16").Value = 9.0879999999999992
$ws.Range("F16").Value = 137.512
$ws.Range("G16").Value = 18.077999999999999
$ws.Range("H16").Value = 2.5419999999999998
$ws.Range("I16").Value = 1.2869999999999999
$ws.Range("J16").Value = 0.31462535379878143
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 37.512
$ws.Range("M16").Value = 105.456
$ws.Range("N16").Value = 2.6669999999999998
$ws.Range("O16").Value = 1.341
$ws.Range("C17").Value = 0.31725840376789194
$ws.Range("D17").Value = 0.99873566969044347
$ws.Range("E17").Value = 9.1150000000000002
$ws.Range("F17").Value = 137.93299999999999
$ws.Range("G17").Value = 18.222000000000001
$ws.Range("H17").Value = 4.5279999999999996
$ws.Range("I17").Value = 2.06
$ws.Range("J17").Value = 0.32071495127050625
$ws.Range("K17").Value = 0.99961597542242697
$ws.Range("L17").Value = 37.933
$ws.Range("M17").Value = 108.976
$ws.Range("N17").Value = 4.3579999999999997
$ws.Range("O17").Value = 2.1760000000000002